# Autogenerated on Sun Feb 01 2015 22:24:41 GMT-0500 (Eastern Standard Time)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Data" to "Summary"
$ws.Name = "Summary"

# Push the existing table (old rows 5-11) down by six rows so it now
# occupies rows 11-17, opening up space for the new "Source Type" line
# at row 9 and leaving rows 18-22 free before the new citation block.
$ws.Rows("5:10").Insert()

# New sub-heading above the table (bold + underlined)
$c = $ws.Range("A9")
$c.Value = "Source Type: SME Associations (Most Widely Used)"
$c.Font.Bold = $true
$c.Font.Underline = $true

# New data point: MSMEs employment (% of total) = 86.6
$c = $ws.Range("D14")
$c.Value = "'86.6"

# New citation block at the bottom of the sheet
$c = $ws.Range("A23")
$c.Value = "CIEN"
$c.Font.Bold = $true

$c = $ws.Range("A24")
$c.Value = "CENTRO DE INVESTIGACIONES ECONÓMICAS NACIONALES, MICRO, PEQUENAS Y MEDIANAS EMPRESAS EN GUATEMALA. Available at http://www.mejoremosguate.org/cms/content/files/diagnosticos/economicos/Lineamientos_PYMES_05-05-2011.pdf"
$c.Font.Italic = $true
